$wb = $excel.ActiveWorkbook

# Target sheet: "SoCDTtiNTY-frgt" (freight) - align with AVL values
$ws = $wb.Worksheets.Item("SoCDTtiNTY-frgt")

# Row 2 (LDVs / HDVs line for freight trucks) -> 1/17
$ws.Range("B2").Formula = "=1/17"
$ws.Range("C2:H2").Formula = "=1/17"

# Row 3 -> 1/19
$ws.Range("B3").Formula = "=1/19"
$ws.Range("C3:H3").Formula = "=1/19"

# Update selection on this sheet to B2:H3 with active cell B2
[void]$ws.Range("B2:H3").Select()

# Make "About" sheet the active/selected tab
[void]$wb.Worksheets.Item("About").Activate()
